$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = ">=20"
$ws.Range("B3").Value = ">=30"
$ws.Range("B4").Value = ">=40"
$ws.Range("B5").Value = ">=50"
$ws.Range("B6").Value = ">=60"
